# cl_neg_with_std.xlsx: renamed the "neg" category labels (markers) in
# column D to their updated phrasing, per "swapped colors and markers,
# updated excel".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Old category label -> new category label.
$map = @{
    "statement neg"     = "make statement neg"
    "cooperation neg"   = "cooperate neg"
    "retreat neg"       = "yield neg"
    "investigation neg" = "investigate neg"
    "demand neg"        = "demand neg"
    "dissaproval neg"   = "disapprove neg"
    "rejection neg"     = "reject neg"
    "threat neg"        = "threaten neg"
    "protest neg"       = "protest neg"
    "force neg"         = "exhibit force neg"
    "relation neg"      = "reduce relations neg"
    "coercion neg"      = "coerce neg"
    "assault neg"       = "assault neg"
    "fight neg"         = "fight neg"
    "hybrid attack neg" = "mass violence neg"
}

for ($r = 2; $r -le 76; $r++) {
    $cur = $ws.Cells.Item($r, 4).Value2
    if ($null -ne $cur -and $map.ContainsKey($cur)) {
        $ws.Cells.Item($r, 4).Value = $map[$cur]
    }
}

# Column D widened (bestFit) to fit the longest new label ("reduce relations neg").
$ws.Columns.Item(4).AutoFit() | Out-Null
$ws.Columns.Item(4).ColumnWidth = 18.59  # best achievable match for stored width 19.42578125

# Selection/viewport moved up the sheet after the edit.
$ws.Range("I46").Select()
